# Daily attendance processing - 2025-11-06 05:22:59
# Normalize the ordering of the "Recorded By" (column G) contributor lists
# for each attendance row: the values scraped from the system were
# reshuffled (old first entry moved to the end of the list) for a known
# set of recurring combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact old-value -> new-value mapping observed for the "Recorded By" column.
$map = @{
    "backup@backdoor.com, system, System" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
